$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H15").Value = 1084.4429
$ws.Range("I15").Value = 1084.4429
$ws.Range("K15").Value = 3253.3287
$ws.Range("M15").Value = -3084.3287
$ws.Range("H33").Value = 183.04347
$ws.Range("I33").Value = 162.17647
$ws.Range("J33").Value = 242.16667
$ws.Range("K33").Value = 162.17647
$ws.Range("L33").Value = 242.16667
$ws.Range("M33").Value = 66.82353000000001
$ws.Range("N33").Value = -700.1666700000001
$ws.Range("H53").Value = 2610.5293
$ws.Range("I53").Value = 2757.4167
$ws.Range("K53").Value = 2757.4167
$ws.Range("M53").Value = -2120.4167
$ws.Range("H76").Value = 100003440
$ws.Range("I76").Value = 250003180
$ws.Range("J76").Value = 3601.6667
$ws.Range("K76").Value = 250003180
$ws.Range("L76").Value = 3601.6667
$ws.Range("M76").Value = -250002865
$ws.Range("N76").Value = -4231.6667
$ws.Range("H79").Value = 100003440
$ws.Range("I79").Value = 250003180
$ws.Range("J79").Value = 3601.6667
$ws.Range("K79").Value = 250003180
$ws.Range("L79").Value = 3601.6667
$ws.Range("M79").Value = -250002088
$ws.Range("N79").Value = -5785.6667
$ws.Range("H86").Value = 5249.5
$ws.Range("J86").Value = 4566.3335
$ws.Range("L86").Value = 4566.3335
$ws.Range("N86").Value = -6812.3335
$ws.Range("H89").Value = 5249.5
$ws.Range("J89").Value = 4566.3335
$ws.Range("L89").Value = 22831.6675
$ws.Range("N89").Value = -34063.6675
$ws.Range("H92").Value = 5499
$ws.Range("I92").Value = 5000
$ws.Range("J92").Value = 5998
$ws.Range("K92").Value = 5000
$ws.Range("L92").Value = 5998
$ws.Range("M92").Value = -3752
$ws.Range("N92").Value = -8494
$ws.Range("H116").Value = 5562.375
$ws.Range("I116").Value = 5951
$ws.Range("K116").Value = 5951
$ws.Range("M116").Value = -2509
$ws.Range("H137").Value = 2208.125
$ws.Range("I137").Value = 2072.6
$ws.Range("J137").Value = 2885.75
$ws.Range("K137").Value = 6217.799999999999
$ws.Range("L137").Value = 8657.25
$ws.Range("M137").Value = -3667.799999999999
$ws.Range("N137").Value = -13757.25
$ws.Range("H138").Value = 2832.2466
$ws.Range("I138").Value = 4628.3335
$ws.Range("J138").Value = 2244.4363
$ws.Range("K138").Value = 13885.0005
$ws.Range("L138").Value = 6733.3089
$ws.Range("M138").Value = -8745.000499999998
$ws.Range("N138").Value = -17013.3089
$ws.Range("H141").Value = 11524.75
$ws.Range("I141").Value = 5015.8
$ws.Range("J141").Value = 16174
$ws.Range("K141").Value = 15047.4
$ws.Range("L141").Value = 48522
$ws.Range("M141").Value = -9867.400000000001
$ws.Range("N141").Value = -58882

# ---- ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H45").Value = 2498.8235
$ws.Range("I45").Value = 2069.4
$ws.Range("K45").Value = 2069.4
$ws.Range("M45").Value = -1692.4
$ws.Range("H88").Value = 9802.833000000001
$ws.Range("I88").Value = 1427.8
$ws.Range("K88").Value = 1427.8
$ws.Range("M88").Value = -1021.8
$ws.Range("H91").Value = 9802.833000000001
$ws.Range("I91").Value = 1427.8
$ws.Range("K91").Value = 1427.8
$ws.Range("M91").Value = -23.79999999999995
$ws.Range("H97").Value = 2867.3333
$ws.Range("I97").Value = 1282.7333
$ws.Range("J97").Value = 6828.8335
$ws.Range("K97").Value = 1282.7333
$ws.Range("L97").Value = 6828.8335
$ws.Range("M97").Value = -786.7333000000001
$ws.Range("N97").Value = -7820.8335
$ws.Range("H102").Value = 3817.7144
$ws.Range("I102").Value = 3785.5
$ws.Range("K102").Value = 3785.5
$ws.Range("M102").Value = -2163.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 1065.6545
$ws.Range("I132").Value = 839.4286
$ws.Range("J132").Value = 2913.1667
$ws.Range("K132").Value = 2518.2858
$ws.Range("L132").Value = 8739.500100000001
$ws.Range("M132").Value = 11.71420000000035
$ws.Range("N132").Value = -13799.5001

# ---- BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H59").Value = 69639
$ws.Range("I59").Value = 69639
$ws.Range("K59").Value = 69639
$ws.Range("M59").Value = -68792
$ws.Range("H86").Value = 9738.5625
$ws.Range("I86").Value = 3272.6428
$ws.Range("K86").Value = 3272.6428
$ws.Range("M86").Value = -2149.6428
$ws.Range("H89").Value = 9738.5625
$ws.Range("I89").Value = 3272.6428
$ws.Range("K89").Value = 16363.214
$ws.Range("M89").Value = -10747.214
$ws.Range("H94").Value = 11014.529
$ws.Range("I94").Value = 6096.1
$ws.Range("J94").Value = 18040.857
$ws.Range("K94").Value = 6096.1
$ws.Range("L94").Value = 18040.857
$ws.Range("M94").Value = -5645.1
$ws.Range("N94").Value = -18942.857
$ws.Range("H105").Value = 3171.2856
$ws.Range("J105").Value = 3633.3333
$ws.Range("L105").Value = 3633.3333
$ws.Range("N105").Value = -7127.3333

# ---- CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H55").Value = 6081
$ws.Range("J55").Value = 6081
$ws.Range("L55").Value = 6081
$ws.Range("M55").Value = -6711
$ws.Range("H105").Value = 1496.5
$ws.Range("I105").Value = 1413
$ws.Range("K105").Value = 1413
$ws.Range("M105").Value = 334
$ws.Range("H140").Value = 130755
$ws.Range("J140").Value = 130755
$ws.Range("L140").Value = 130755
$ws.Range("N140").Value = -141115

# ---- GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H7").Value = 8343483
$ws.Range("I7").Value = 12513225
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 12513225
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -12513113
$ws.Range("N7").Value = -4224
$ws.Range("H8").Value = 8343483
$ws.Range("I8").Value = 12513225
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 12513225
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = -12513086
$ws.Range("N8").Value = -4278
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H111").Value = 145332.67
$ws.Range("J111").Value = 145332.67
$ws.Range("L111").Value = 145332.67
$ws.Range("N111").Value = -151466.67
$ws.Range("H113").Value = 4992
$ws.Range("I113").Value = 4634.8335
$ws.Range("J113").Value = 5170.5835
$ws.Range("K113").Value = 4634.8335
$ws.Range("L113").Value = 5170.5835
$ws.Range("M113").Value = -2464.8335
$ws.Range("N113").Value = -9510.583500000001

# ---- LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 2750
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2750
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2974
$ws.Range("H13").Value = 9000
$ws.Range("I13").Value = 9000
$ws.Range("K13").Value = 9000
$ws.Range("M13").Value = -8860
$ws.Range("H22").Value = 1368.7778
$ws.Range("I22").Value = 887.7143
$ws.Range("J22").Value = 1537.15
$ws.Range("K22").Value = 887.7143
$ws.Range("L22").Value = 1537.15
$ws.Range("M22").Value = -592.7143
$ws.Range("N22").Value = -2127.15
$ws.Range("H27").Value = 1368.7778
$ws.Range("I27").Value = 887.7143
$ws.Range("J27").Value = 1537.15
$ws.Range("K27").Value = 887.7143
$ws.Range("L27").Value = 1537.15
$ws.Range("M27").Value = -780.7143
$ws.Range("N27").Value = -1751.15
$ws.Range("H29").Value = 38500
$ws.Range("I29").Value = 38500
$ws.Range("K29").Value = 38500
$ws.Range("M29").Value = -38205
$ws.Range("H59").Value = 32950
$ws.Range("J59").Value = 32950
$ws.Range("L59").Value = 32950
$ws.Range("N59").Value = -34258
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8250
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 2299.625
$ws.Range("J132").Value = 2647.1
$ws.Range("L132").Value = 7941.299999999999
$ws.Range("N132").Value = -13001.3
$ws.Range("H136").Value = 2082.8914
$ws.Range("I136").Value = 1961.9333
$ws.Range("J136").Value = 2309.6875
$ws.Range("K136").Value = 5885.7999
$ws.Range("L136").Value = 6929.0625
$ws.Range("M136").Value = -3335.7999
$ws.Range("N136").Value = -12029.0625

# ---- WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5576
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H109").Value = 76376.664
$ws.Range("J109").Value = 76376.664
$ws.Range("L109").Value = 76376.664
$ws.Range("N109").Value = -79150.664
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 4082.2166
$ws.Range("I132").Value = 3784.5715
$ws.Range("J132").Value = 4498.92
$ws.Range("K132").Value = 11353.7145
$ws.Range("L132").Value = 13496.76
$ws.Range("M132").Value = -8823.7145
$ws.Range("N132").Value = -18556.76
$ws.Range("H136").Value = 5093.579
$ws.Range("I136").Value = 9866.333000000001
$ws.Range("J136").Value = 2890.7693
$ws.Range("K136").Value = 29598.999
$ws.Range("L136").Value = 8672.3079
$ws.Range("M136").Value = -27048.999
$ws.Range("N136").Value = -13772.3079
